$wb = $excel.ActiveWorkbook

# The existing "Greece" sheet is the structural template (column widths,
# row heights, dimension) that the new per-market sheets are based on.
$template = $wb.Worksheets.Item("Greece")

# --- Norway -----------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $lastSheet)
$norway = $wb.Worksheets.Item($wb.Worksheets.Count)
$norway.Name = "Norway"
$norway.Range("B4").Value = "NGC-2931/T3058"
$norway.Range("B2").Value = "Norway Market"
$norway.Cells.Select()

# --- Poland -------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $lastSheet)
$poland = $wb.Worksheets.Item($wb.Worksheets.Count)
$poland.Name = "Poland"
$poland.Range("B4").Value = "NGC-2920/T3101"
$poland.Range("B2").Value = "Poland Market"
$poland.Cells.Select()

$norway.Activate()
